$wb = $excel.ActiveWorkbook

# --- Sheet "isa_template": bump template Version ---
$ws1 = $wb.Worksheets.Item("isa_template")
$ws1.Range("B4").Value = "1.0.3"

# --- Sheet "New Table": rename building block headers ---
$ws2 = $wb.Worksheets.Item("New Table")
$ws2.Range("E1").Value = "Component [NMR instrument]"
$ws2.Range("H1").Value = "Component [NMR probe]"
$ws2.Range("U1").Value = "Output [Data]"

# --- Update term accession number URLs to bioregistry.io style ---
$ws2.Range("D2").Value = "https://bioregistry.io/OBI:0000623"
$ws2.Range("G2").Value = "https://bioregistry.io/OBI:0000558"
$ws2.Range("T2").Value = "https://bioregistry.io/UO:0000228"
